$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Species concentrations reverted to 1e-9 (row 2-7, column C) ---
$ws.Range("C2").Value = 0.000000001
$ws.Range("C3").Value = 0.000000001
$ws.Range("C4").Value = 0.000000001
$ws.Range("C5").Value = 0.000000001
$ws.Range("C6").Value = 0.000000001
$ws.Range("C7").Value = 0.000000001

# --- Receptor densities (row 8-12) ---
$ws.Range("C8").Formula = "=1600"
$ws.Range("C9").Formula = "=4100"
$ws.Range("C10").Value = 44100
$ws.Range("C11").Value = 5000
$ws.Range("C12").Value = 5300

# --- Association/dissociation rates (rows 32-37) ---
$ws.Range("C32").Value = 11000000
$ws.Range("C34").Value = 4600000
$ws.Range("C36").Value = 1300000
$ws.Range("C37").Value = 0.0035

# --- Row 43 ---
$ws.Range("C43").Value = 0.0000915

# --- Row 55 ---
$ws.Range("C55").Value = 0.000304

# --- Rows 64-65: replace formulas with plain values ---
$ws.Range("C64").Value = 331000
$ws.Range("C65").Value = 0.000755

# --- Restore sheet view/selection state ---
$null = $ws.Range("C5").Select()
